$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 43.142857
$ws.Range("I5").Value = 45.333332
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 45.333332
$ws.Range("L5").Value = 30
$ws.Range("M5").Value = 69.666668
$ws.Range("N5").Value = -260

$ws.Range("H32").Value = 844.61536
$ws.Range("J32").Value = 997.5
$ws.Range("L32").Value = 997.5
$ws.Range("N32").Value = -1649.5

$ws.Range("H40").Value = 3573.2
$ws.Range("I40").Value = 4247.75
$ws.Range("J40").Value = 875
$ws.Range("K40").Value = 4247.75
$ws.Range("L40").Value = 875
$ws.Range("M40").Value = -4072.75
$ws.Range("N40").Value = -1225

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0

$ws.Range("H137").Value = 2471.2144
$ws.Range("I137").Value = 1876.5
$ws.Range("J137").Value = 3264.1667
$ws.Range("K137").Value = 5629.5
$ws.Range("L137").Value = 9792.500100000001
$ws.Range("M137").Value = -3079.5
$ws.Range("N137").Value = -14892.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 23634.455
$ws.Range("J37").Value = 24997.666
$ws.Range("L37").Value = 24997.666
$ws.Range("N37").Value = -25543.666

$ws.Range("H55").Value = 20000

$ws.Range("H63").Value = 2361.8572
$ws.Range("I63").Value = 1706.6
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 1706.6
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -1020.6
$ws.Range("N63").Value = -5372

$ws.Range("H66").Value = 2361.8572
$ws.Range("I66").Value = 1706.6
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 8533
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -5101
$ws.Range("N66").Value = -26864

$ws.Range("H110").Value = 2525.1667
$ws.Range("I110").Value = 950.2727
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 950.2727
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = 1094.7273
$ws.Range("N110").Value = -9090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1360.3334
$ws.Range("I20").Value = 1177.3334
$ws.Range("K20").Value = 1177.3334
$ws.Range("M20").Value = -930.3334

$ws.Range("H35").Value = 14997
$ws.Range("J35").Value = 14997
$ws.Range("L35").Value = 14997
$ws.Range("N35").Value = -15617

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 15589.546
$ws.Range("J41").Value = 20712.143
$ws.Range("L41").Value = 20712.143
$ws.Range("N41").Value = -21568.143

$ws.Range("H50").Value = 29982.908
$ws.Range("J50").Value = 29982.908
$ws.Range("L50").Value = 29982.908
$ws.Range("N50").Value = -31232.908

$ws.Range("H59").Value = 30925.555
$ws.Range("I59").Value = 20381.4
$ws.Range("J59").Value = 34981
$ws.Range("K59").Value = 20381.4
$ws.Range("L59").Value = 34981
$ws.Range("M59").Value = -19236.4
$ws.Range("N59").Value = -37271

$ws.Range("H60").Value = 21638.46
$ws.Range("I60").Value = 3243
$ws.Range("J60").Value = 24983.092
$ws.Range("K60").Value = 3243
$ws.Range("L60").Value = 24983.092
$ws.Range("M60").Value = -2732
$ws.Range("N60").Value = -26005.092

$ws.Range("H62").Value = 9250
$ws.Range("I62").Value = 9250
$ws.Range("K62").Value = 9250
$ws.Range("M62").Value = -8626

$ws.Range("H65").Value = 9250
$ws.Range("I65").Value = 9250
$ws.Range("K65").Value = 46250
$ws.Range("M65").Value = -43130

$ws.Range("H68").Value = 39298.715
$ws.Range("J68").Value = 39993.383
$ws.Range("L68").Value = 39993.383
$ws.Range("N68").Value = -41491.383

$ws.Range("H71").Value = 39298.715
$ws.Range("J71").Value = 39993.383
$ws.Range("L71").Value = 119980.149
$ws.Range("N71").Value = -127468.149

$ws.Range("H74").Value = 38191.438
$ws.Range("J74").Value = 38191.438
$ws.Range("L74").Value = 38191.438
$ws.Range("N74").Value = -39939.438

$ws.Range("H77").Value = 38191.438
$ws.Range("J77").Value = 38191.438
$ws.Range("L77").Value = 114574.314
$ws.Range("N77").Value = -123310.314

$ws.Range("H86").Value = 4648.6
$ws.Range("I86").Value = 4609.5557
$ws.Range("K86").Value = 4609.5557
$ws.Range("M86").Value = -3486.5557

$ws.Range("H89").Value = 4648.6
$ws.Range("I89").Value = 4609.5557
$ws.Range("K89").Value = 23047.7785
$ws.Range("M89").Value = -17431.7785

$ws.Range("H103").Value = 7075
$ws.Range("I103").Value = 7075
$ws.Range("K103").Value = 7075
$ws.Range("M103").Value = -5903

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2220.1875
$ws.Range("I131").Value = 1302.875
$ws.Range("J131").Value = 3137.5
$ws.Range("K131").Value = 3908.625
$ws.Range("L131").Value = 9412.5
$ws.Range("M131").Value = 1131.375
$ws.Range("N131").Value = -19492.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 19070.637
$ws.Range("I43").Value = 8000
$ws.Range("K43").Value = 8000
$ws.Range("M43").Value = -7849

$ws.Range("H97").Value = 894.93335
$ws.Range("I97").Value = 894.93335
$ws.Range("K97").Value = 894.93335
$ws.Range("M97").Value = -398.93335

$ws.Range("H102").Value = 1280.9375
$ws.Range("I102").Value = 1280.9375
$ws.Range("K102").Value = 1280.9375
$ws.Range("M102").Value = 341.0625

$ws.Range("H107").Value = 162.8
$ws.Range("I107").Value = 162.8
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 162.8
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 1757.2

$ws.Range("H122").Value = 9639299
$ws.Range("I122").Value = 11389044
$ws.Range("J122").Value = 15699
$ws.Range("K122").Value = 34167132
$ws.Range("L122").Value = 47097
$ws.Range("M122").Value = -34164682
$ws.Range("N122").Value = -51997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 2500
$ws.Range("I42").Value = 2500
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 2500
$ws.Range("L42").ClearContents()
$ws.Range("N42").Value = 0
$ws.Range("M42").Value = -1937

$ws.Range("H49").Value = 2500
$ws.Range("I49").Value = 2500
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 2500
$ws.Range("L49").ClearContents()
$ws.Range("N49").Value = 0
$ws.Range("M49").Value = -2353

$ws.Range("H68").Value = 3100.4
$ws.Range("I68").Value = 2667.3333
$ws.Range("J68").Value = 3750
$ws.Range("K68").Value = 2667.3333
$ws.Range("L68").Value = 3750
$ws.Range("M68").Value = -1918.3333
$ws.Range("N68").Value = -5248

$ws.Range("H71").Value = 3100.4
$ws.Range("I71").Value = 2667.3333
$ws.Range("J71").Value = 3750
$ws.Range("K71").Value = 13336.6665
$ws.Range("L71").Value = 18750
$ws.Range("M71").Value = -9592.666499999999
$ws.Range("N71").Value = -26238

$ws.Range("H93").Value = 839.8
$ws.Range("I93").Value = 799.5
$ws.Range("K93").Value = 799.5
$ws.Range("M93").Value = 448.5

$ws.Range("H100").Value = 2854.7144
$ws.Range("I100").Value = 2854.7144
$ws.Range("K100").Value = 2854.7144
$ws.Range("M100").Value = -2313.7144

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 40000
$ws.Range("J95").Value = 40000
$ws.Range("L95").Value = 40000
$ws.Range("N95").Value = -45492

$ws.Range("H96").Value = 1499.3334
$ws.Range("I96").Value = 1251.5
$ws.Range("K96").Value = 1251.5
$ws.Range("M96").Value = 121.5
